$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.413.90"
$ws.Range("E2").Value = "  +1.04%  "
$ws.Range("D3").Value = "1.637.94"
$ws.Range("E3").Value = "  +2.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.36"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3745"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.95%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.02"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3632"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.258"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08134"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.90"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.619"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001273"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.336"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.04%  "
$ws.Range("D17").Value = "1.630.79"
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.43"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06905"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.529"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "23.417.28"
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.75"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.097"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.420"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.21"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.03"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.332"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.48"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.303"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.73%  "
$ws.Range("D32").Value = "1.812.56"
$ws.Range("E32").Value = "  +1.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.744"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9587"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02850"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.29"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.07295"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2524"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08809"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.105"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.49%  "
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7070"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.44"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.16"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6538"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.327"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.0000"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.012"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07966"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.60"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.70%  "
$ws.Range("E51").Value = "  +0.22%  "
